# Auto-generated edit script applying cell updates per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "28.016.90"
Set-TextValue $ws "E2" "  -1.97%  "
Set-TextValue $ws "D3" "1.830.63"
Set-TextValue $ws "E3" "  -1.10%  "
Set-TextValue $ws "E4" "  -0.11%  "
Set-TextValue $ws "D5" "325.22"
Set-TextValue $ws "E5" "  -3.03%  "
Set-TextValue $ws "E6" "  -0.10%  "
Set-TextValue $ws "D7" "0.4644"
Set-TextValue $ws "E7" "  -0.31%  "
Set-TextValue $ws "D8" "0.3876"
Set-TextValue $ws "E8" "  -1.09%  "
Set-TextValue $ws "D9" "0.07870"
Set-TextValue $ws "E9" "  -0.27%  "
Set-TextValue $ws "D10" "0.9598"
Set-TextValue $ws "E10" "  -2.48%  "
Set-TextValue $ws "E11" "  -1.55%  "
Set-TextValue $ws "D12" "1.801.38"
Set-TextValue $ws "E12" "  -6.95%  "
Set-TextValue $ws "D13" "5.675"
Set-TextValue $ws "E13" "  -2.94%  "
Set-TextValue $ws "D14" "6.906"
Set-TextValue $ws "E14" "  -1.58%  "
Set-TextValue $ws "D15" "0.06779"
Set-TextValue $ws "E15" "  -0.60%  "
Set-TextValue $ws "D16" "87.27"
Set-TextValue $ws "E16" "  -0.43%  "
Set-TextValue $ws "E17" "  -0.15%  "
Set-TextValue $ws "D18" "0.000009923"
Set-TextValue $ws "E18" "  -1.87%  "
Set-TextValue $ws "D19" "16.63"
Set-TextValue $ws "E19" "  -2.32%  "
Set-TextValue $ws "D20" "1.002"
Set-TextValue $ws "D21" "28.026.74"
Set-TextValue $ws "D22" "5.316"
Set-TextValue $ws "E22" "  -1.64%  "
Set-TextValue $ws "D23" "10.98"
Set-TextValue $ws "E23" "  -2.42%  "
Set-TextValue $ws "D24" "2.095"
Set-TextValue $ws "E24" "  -1.42%  "
Set-TextValue $ws "D25" "2.077.78"
Set-TextValue $ws "E25" "  -3.51%  "
Set-TextValue $ws "D26" "153.84"
Set-TextValue $ws "E26" "  +0.24%  "
Set-TextValue $ws "D27" "19.13"
Set-TextValue $ws "E27" "  -1.25%  "
Set-TextValue $ws "D28" "5.738"
Set-TextValue $ws "E28" "  -7.15%  "
Set-TextValue $ws "D29" "1.974"
Set-TextValue $ws "E29" "  -2.35%  "
Set-TextValue $ws "D30" "117.34"
Set-TextValue $ws "E30" "  -0.09%  "
Set-TextValue $ws "D31" "0.9360"
Set-TextValue $ws "E31" "  -4.13%  "
Set-TextValue $ws "D32" "0.09262"
Set-TextValue $ws "E32" "  -1.90%  "
Set-TextValue $ws "E33" "  -1.51%  "
Set-TextValue $ws "E34" "  -2.34%  "
Set-TextValue $ws "D35" "3.287"
Set-TextValue $ws "E35" "  -6.15%  "
Set-TextValue $ws "D36" "0.05862"
Set-TextValue $ws "E36" "  -4.34%  "
Set-TextValue $ws "D37" "0.02143"
Set-TextValue $ws "E37" "  -2.34%  "
Set-TextValue $ws "D38" "1.145"
Set-TextValue $ws "E38" "  -1.53%  "
Set-TextValue $ws "D39" "7.788"
Set-TextValue $ws "E39" "  +2.60%  "
Set-TextValue $ws "D40" "0.5589"
Set-TextValue $ws "E40" "  -1.99%  "
Set-TextValue $ws "D41" "9.862"
Set-TextValue $ws "E41" "  -2.50%  "
Set-TextValue $ws "D42" "0.1763"
Set-TextValue $ws "E42" "  -1.52%  "
Set-TextValue $ws "E43" "  -1.53%  "
Set-TextValue $ws "D44" "0.5263"
Set-TextValue $ws "E44" "  -2.29%  "
Set-TextValue $ws "D45" "0.07000"
Set-TextValue $ws "E45" "  -2.11%  "
Set-TextValue $ws "D46" "2.118"
Set-TextValue $ws "E46" "  -11.00%  "
Set-TextValue $ws "B47" "WEMIXToken"
Set-TextValue $ws "C47" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D47" "1.119"
Set-TextValue $ws "E47" "  -11.09%  "
Set-TextValue $ws "B48" "NEARProtocol"
Set-TextValue $ws "C48" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D48" "1.832"
Set-TextValue $ws "E48" "  -3.91%  "
Set-TextValue $ws "B49" "Quant"
Set-TextValue $ws "C49" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws "D49" "113.04"
Set-TextValue $ws "E49" "  -0.45%  "
Set-TextValue $ws "E50" "  -0.12%  "
Set-TextValue $ws "E51" "  +0.41%  "

Write-Host "Applied 95 cell updates"
